# Apply the changes described by the commit:
# "Changes to code to run all test without fail expect for leave test case"
#
# 1. Admin sheet: cell D2 gets a new generated id/key value.
# 2. Jira sheet: cell B2 gets a new error-description hash value.
# 3. Jira sheet: cell A3 changes from the previous run-key to "Leave_Accept".
# 4. Jira sheet: a new row (row 4) is appended recording "PersonalDetails".

$wb = $excel.ActiveWorkbook

$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "7064986A"

$wsJira = $wb.Worksheets.Item("Jira")
$wsJira.Range("B2").Value = "0E8CAC993B05BCDF8711C30890361286"
$wsJira.Range("A3").Value = "Leave_Accept"
$wsJira.Range("A4").Value = "PersonalDetails"

# Give the new B4 cell the same (default) formatting as the existing blank
# cell above it (B3) so the row is materialized with the expected style.
$wsJira.Range("B3").Copy()
$wsJira.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host ("Admin!D2 = " + $wsAdmin.Range("D2").Text)
Write-Host ("Jira!B2 = " + $wsJira.Range("B2").Text)
Write-Host ("Jira!A3 = " + $wsJira.Range("A3").Text)
Write-Host ("Jira!A4 = " + $wsJira.Range("A4").Text)
Write-Host ("Jira dimension used range rows = " + $wsJira.UsedRange.Rows.Count)
